# Updates the cryptos list (Price / Volume(1h) columns, plus a row swap
# between TheGraph and Kaspa) to reflect the latest scrape.
#
# NOTE: Column D "Price" values are text that often look numeric
# (e.g. "68.191.01", "0.999", "0.0000353"). Assigning such a string
# straight to .Value would make Excel auto-convert it to a real number
# and silently mangle it (drop trailing zeros, switch to scientific
# notation, merge thousand-separator dots, etc). Prefixing the literal
# with a leading apostrophe forces Excel to keep it as text, exactly
# like a user typing `'68.191.01` into a cell - the apostrophe itself
# is not stored as part of the value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.191.01"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "'3.887.92"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'482.71"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "'145.60"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +0.72%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.742"
$ws.Range("E9").Value = "  +3.21%  "

$ws.Range("D10").Value = "'0.180"
$ws.Range("E10").Value = "  +7.82%  "

$ws.Range("D11").Value = "'0.0000353"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").Value = "'43.09"
$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("D13").Value = "'10.50"
$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").Value = "'4.504.44"
$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "'3.897.95"
$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").Value = "'14.24"
$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("E18").Value = "  +1.64%  "

$ws.Range("E19").Value = "  +0.81%  "

$ws.Range("D20").Value = "'68.219.56"
$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("D21").Value = "'428.90"
$ws.Range("E21").Value = "  -0.44%  "

$ws.Range("E22").Value = "  +8.24%  "

$ws.Range("D23").Value = "'14.76"
$ws.Range("E23").Value = "  +1.98%  "

$ws.Range("D24").Value = "'12.34"
$ws.Range("E24").Value = "  +18.15%  "

$ws.Range("D25").Value = "'88.75"
$ws.Range("E25").Value = "  +2.41%  "

$ws.Range("E26").Value = "  +3.13%  "

$ws.Range("D27").Value = "'11.00"
$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("D28").Value = "'37.17"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("D30").Value = "'718.68"
$ws.Range("E30").Value = "  +0.92%  "

$ws.Range("D31").Value = "'13.47"
$ws.Range("E31").Value = "  +2.16%  "

$ws.Range("E32").Value = "  +0.50%  "

$ws.Range("D33").Value = "'2.92"
$ws.Range("E33").Value = "  +3.02%  "

$ws.Range("D34").Value = "'61.77"
$ws.Range("E34").Value = "  +6.21%  "

$ws.Range("D35").Value = "'0.0₃0881"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").Value = "'6.05"
$ws.Range("E36").Value = "  +10.30%  "

$ws.Range("D37").Value = "'40.80"
$ws.Range("E37").Value = "  -0.94%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.146"
$ws.Range("E38").Value = "  -3.10%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.398"
$ws.Range("E39").Value = "  +17.22%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("E41").Value = "  +6.30%  "

$ws.Range("D42").Value = "'2.99"
$ws.Range("E42").Value = "  +8.41%  "

$ws.Range("E43").Value = "  +3.76%  "

$ws.Range("E44").Value = "  -1.79%  "

$ws.Range("D45").Value = "'0.142"
$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("D46").Value = "'3.37"
$ws.Range("E46").Value = "  +7.70%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "'0.0₆0352"
$ws.Range("E48").Value = "  +28.01%  "

$ws.Range("E49").Value = "  -0.78%  "

$ws.Range("E50").Value = "  -2.02%  "

$ws.Range("D51").Value = "'144.32"
$ws.Range("E51").Value = "  -2.60%  "
